# Lecture 16 Section 1 - 2nd commit
# Adds a new "OpenAccountTest" worksheet, extra customer rows on
# "AddCustomerTest", and restores the original selections/active sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Add the new "OpenAccountTest" worksheet after AddCustomerTest
#    (added/populated first so its strings land first in the shared
#    string table, matching the order they were authored in)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Raman Arora"
$ws2.Range("B2").Value = "Rupee"

# wrap the customer name cell
$ws2.Range("A2").WrapText = $true

# best-fit the first column like Excel's Format > AutoFit Column Width
$ws2.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------
# 2. AddCustomerTest: append three more customer rows (3-5)
# ---------------------------------------------------------------
$ws1.Range("A3").Value = "Rahul"
$ws1.Range("B3").Value = "Arora"
$ws1.Range("C3").Value = "A234wd"
$ws1.Range("D3").Value = "Customer added successfully"

$ws1.Range("A4").Value = "Ishita"
$ws1.Range("B4").Value = "Arora"
$ws1.Range("C4").Value = "A234wd"
$ws1.Range("D4").Value = "Customer added successfully"

$ws1.Range("A5").Value = "Rohit"
$ws1.Range("B5").Value = "Sehgal"
$ws1.Range("C5").Value = "A234wd"
$ws1.Range("D5").Value = "Customer added successfully"

# ---------------------------------------------------------------
# 3. Restore view state: OpenAccountTest keeps selection K13, while
#    AddCustomerTest (the originally active tab) ends up selected at B6.
# ---------------------------------------------------------------
$ws2.Range("K13").Select() | Out-Null

$ws1.Activate()
$ws1.Range("B6").Select() | Out-Null
